$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $wb.Windows.Item(1)
Write-Output $win
Write-Output $win.TabRatio
$win.TabRatio = 0.989
Write-Output $win.TabRatio
